$wb = $excel.ActiveWorkbook

# --- ADI_scaling: new "N = 3840" header label + a second scaling experiment block (N = 7680) ---
$ws = $wb.Worksheets.Item("ADI_scaling")

$ws.Range("B1").Value = "N = 3840"

$ws.Range("A15").Value = "Strong"
$ws.Range("B15").Value = "N = 7680"

$ws.Range("A16").Value = "n_threads"
$ws.Range("B16").Value = "cycles"

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = 1628927133144
$ws.Range("C17").Formula = '=$B$17/B17'
$ws.Range("C17").NumberFormat = "0.00"

$ws.Range("A18").Value = 2
$ws.Range("B18").Value = 875899735704
$ws.Range("C18").Formula = '=$B$17/B18'
$ws.Range("C18").NumberFormat = "0.00"

$ws.Range("A19").Value = 4
$ws.Range("B19").Value = 452355152308
$ws.Range("C19").Formula = '=$B$17/B19'
$ws.Range("C19").NumberFormat = "0.00"

$ws.Range("A20").Value = 6
$ws.Range("B20").Value = 327365800187
$ws.Range("C20").Formula = '=$B$17/B20'
$ws.Range("C20").NumberFormat = "0.00"

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 254374423524
$ws.Range("C21").Formula = '=$B$17/B21'
$ws.Range("C21").NumberFormat = "0.00"

$ws.Range("A22").Value = 10
$ws.Range("B22").Value = 231848527750
$ws.Range("C22").Formula = '=$B$17/B22'
$ws.Range("C22").NumberFormat = "0.00"

$ws.Range("A23").Value = 12
$ws.Range("B23").Value = 220340724398
$ws.Range("C23").Formula = '=$B$17/B23'
$ws.Range("C23").NumberFormat = "0.00"

$ws.Range("A24").Value = 16
$ws.Range("B24").Value = 237683005100
$ws.Range("C24").Formula = '=$B$17/B24'
$ws.Range("C24").NumberFormat = "0.00"

$ws.Range("A25").Value = 20
$ws.Range("B25").Value = 218365048581
$ws.Range("C25").Formula = '=$B$17/B25'
$ws.Range("C25").NumberFormat = "0.00"

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 226347044864
$ws.Range("C26").Formula = '=$B$17/B26'
$ws.Range("C26").NumberFormat = "0.00"

[void]$ws.Activate()
[void]$ws.Range("D26").Select()

# --- RW_scaling: record the N used for this run (3840) next to the header ---
$ws2 = $wb.Worksheets.Item("RW_scaling")
$ws2.Range("B1").Value = 3840

[void]$ws2.Activate()
[void]$ws2.Range("F17").Select()
